# TOD-E norms run, POM rescale, 24 cell demo strat
#
# The raw-score -> standard-score lookup tab "7.0-9.3" is split into four
# narrower age-band tabs: 7.0-7.5, 7.6-7.11, 8.0-8.5, 8.6-9.3. The first
# (widest) tab is rescaled in place; the other three are new tabs, each
# copied from the previous one (so every tab inherits the header row's
# bold/centered style) and inserted immediately to its right.

$wb = $excel.ActiveWorkbook
$sheet70to75 = $wb.Worksheets.Item("7.0-9.3")

$sheet70to75.Copy([System.Reflection.Missing]::Value, $sheet70to75)
$sheet76to711 = $wb.ActiveSheet
$sheet76to711.Name = "7.6-7.11"

$sheet76to711.Copy([System.Reflection.Missing]::Value, $sheet76to711)
$sheet80to85 = $wb.ActiveSheet
$sheet80to85.Name = "8.0-8.5"

$sheet80to85.Copy([System.Reflection.Missing]::Value, $sheet80to85)
$sheet86to93 = $wb.ActiveSheet
$sheet86to93.Name = "8.6-9.3"

# Rescale the original tab in place and rename it to the narrower band.
$sheet70to75.Name = "7.0-7.5"
$values70to75 = @(64,66,68,70,72,74,76,77,79,81,82,84,86,87,89,91,93,94,96,98,100,102,104,106,109,111,114,117,121,125)
for ($i = 0; $i -lt $values70to75.Length; $i++) {
    $sheet70to75.Cells.Item($i + 2, 2).Value = $values70to75[$i]
}

$values76to711 = @(59,62,64,66,67,69,71,73,74,76,78,79,81,83,85,86,88,90,92,94,96,98,100,102,105,108,111,115,121,127)
for ($i = 0; $i -lt $values76to711.Length; $i++) {
    $sheet76to711.Cells.Item($i + 2, 2).Value = $values76to711[$i]
}

$values80to85 = @(55,58,60,62,63,65,67,69,70,72,74,75,77,79,80,82,84,86,88,90,92,94,97,99,102,106,110,116,123,128)
for ($i = 0; $i -lt $values80to85.Length; $i++) {
    $sheet80to85.Cells.Item($i + 2, 2).Value = $values80to85[$i]
}

$values86to93 = @(51,53,55,57,59,61,62,64,66,67,69,71,72,74,76,78,79,81,83,85,88,90,93,96,101,106,113,118,122,124)
for ($i = 0; $i -lt $values86to93.Length; $i++) {
    $sheet86to93.Cells.Item($i + 2, 2).Value = $values86to93[$i]
}
